# Update "want to go" counts (column F) for a few events.
# 展览 (Exhibition) sheet: rows 4, 17, 20
# 全部类型 (All Types) sheet: rows 4, 20, 24

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 8644
$wsExhibition.Range("F17").Value = 6111
$wsExhibition.Range("F20").Value = 2247

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F4").Value = 8644
$wsAllTypes.Range("F20").Value = 6111
$wsAllTypes.Range("F24").Value = 2247
